$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A49").Value = Get-Date -Year 2023 -Month 11 -Day 5 -Hour 0 -Minute 0 -Second 0
$ws.Range("A49").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B49").Value = "11:16"
$ws.Range("C49").Value = 57.5
$ws.Range("D49").Value = "natura"
